$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Neg_Change")
$ws2 = $wb.Worksheets.Item("Pos_Change")

# Sheet1 (Neg_Change) updates
$ws1.Range("A2").Value = "CIPLA"
$ws1.Range("B2").Value = 1483
$ws1.Range("C2").Value = 1493.8
$ws1.Range("D2").Value = 1478.3
$ws1.Range("E2").Value = 1487.3
$ws1.Range("F2").Value = 920270
$ws1.Range("G2").Value = 1882304
$ws1.Range("H2").Value = -0.5110938509401245
$ws1.Range("I2").Value = "CIPLA"
$ws1.Range("A3").Value = "SUNPHARMA"
$ws1.Range("B3").Value = 1590
$ws1.Range("C3").Value = 1602
$ws1.Range("D3").Value = 1580.9
$ws1.Range("E3").Value = 1587.8
$ws1.Range("F3").Value = 1654424
$ws1.Range("G3").Value = 3435674
$ws1.Range("H3").Value = -0.5184572226584944
$ws1.Range("I3").Value = "SUNPHARMA"
$ws1.Range("A4").Value = "HEROMOTOCO"
$ws1.Range("B4").Value = 4660
$ws1.Range("C4").Value = 4712
$ws1.Range("D4").Value = 4588.1
$ws1.Range("E4").Value = 4600
$ws1.Range("F4").Value = 972461
$ws1.Range("G4").Value = 2121194
$ws1.Range("H4").Value = -0.5415501835287108
$ws1.Range("I4").Value = "HEROMOTOCO"
$ws1.Range("A5").Value = "TORNTPHARM"
$ws1.Range("B5").Value = 3600
$ws1.Range("C5").Value = 3619.2
$ws1.Range("D5").Value = 3562.8
$ws1.Range("E5").Value = 3578.9
$ws1.Range("F5").Value = 178509
$ws1.Range("G5").Value = 402653
$ws1.Range("H5").Value = -0.5566679001522402
$ws1.Range("I5").Value = "TORNTPHARM"
$ws1.Range("A6").Value = "RECLTD"
$ws1.Range("B6").Value = 386.5
$ws1.Range("C6").Value = 387
$ws1.Range("D6").Value = 380.1
$ws1.Range("E6").Value = 381.4
$ws1.Range("F6").Value = 3841685
$ws1.Range("G6").Value = 8126214
$ws1.Range("H6").Value = -0.5272478672109792
$ws1.Range("I6").Value = "RECLTD"
$ws1.Range("A7").Value = "TVSMOTOR"
$ws1.Range("B7").Value = 3014.9
$ws1.Range("C7").Value = 3024.2
$ws1.Range("D7").Value = 2960
$ws1.Range("E7").Value = 2978.8
$ws1.Range("F7").Value = 387198
$ws1.Range("G7").Value = 830216
$ws1.Range("H7").Value = -0.5336177573065323
$ws1.Range("I7").Value = "TVSMOTOR"
$ws1.Range("A8").Value = "DLF"
$ws1.Range("B8").Value = 759.2
$ws1.Range("C8").Value = 768.6
$ws1.Range("D8").Value = 743.15
$ws1.Range("E8").Value = 745.2
$ws1.Range("F8").Value = 1923910
$ws1.Range("G8").Value = 3898418
$ws1.Range("H8").Value = -0.5064895555068748
$ws1.Range("I8").Value = "DLF"
$ws1.Range("A9").Value = "MOTHERSON"
$ws1.Range("B9").Value = 94.12
$ws1.Range("C9").Value = 94.54000000000001
$ws1.Range("D9").Value = 91
$ws1.Range("E9").Value = 91.56
$ws1.Range("F9").Value = 9884650
$ws1.Range("G9").Value = 24236301
$ws1.Range("H9").Value = -0.592155172524058
$ws1.Range("I9").Value = "MOTHERSON"
$ws1.Range("A10").Value = "APLAPOLLO"
$ws1.Range("B10").Value = 1598
$ws1.Range("C10").Value = 1607.5
$ws1.Range("D10").Value = 1575
$ws1.Range("E10").Value = 1577.7
$ws1.Range("F10").Value = 257517
$ws1.Range("G10").Value = 630943
$ws1.Range("H10").Value = -0.5918537807694197
$ws1.Range("I10").Value = "APLAPOLLO"
$ws1.Range("A11").Value = "IRCTC"
$ws1.Range("B11").Value = 727.4
$ws1.Range("C11").Value = 730.9
$ws1.Range("D11").Value = 717
$ws1.Range("E11").Value = 717.8
$ws1.Range("F11").Value = 359497
$ws1.Range("G11").Value = 804696
$ws1.Range("H11").Value = -0.5532511656575899
$ws1.Range("I11").Value = "IRCTC"
$ws1.Range("A12").Value = "BANKINDIA"
$ws1.Range("B12").Value = 112
$ws1.Range("C12").Value = 112.05
$ws1.Range("D12").Value = 110.18
$ws1.Range("E12").Value = 110.4
$ws1.Range("F12").Value = 3578107
$ws1.Range("G12").Value = 8568226
$ws1.Range("H12").Value = -0.5823981533633683
$ws1.Range("I12").Value = "BANKINDIA"
$ws1.Range("A13").Value = "BLUESTARCO"
$ws1.Range("B13").Value = 1830
$ws1.Range("C13").Value = 1830.3
$ws1.Range("D13").Value = 1763
$ws1.Range("E13").Value = 1802.7
$ws1.Range("F13").Value = 631672
$ws1.Range("G13").Value = 1574061
$ws1.Range("H13").Value = -0.5986991609600899
$ws1.Range("I13").Value = "BLUESTARCO"
$ws1.Range("A14").Value = "IGL"
$ws1.Range("B14").Value = 204.9
$ws1.Range("C14").Value = 206.26
$ws1.Range("D14").Value = 201.8
$ws1.Range("E14").Value = 202.5
$ws1.Range("F14").Value = 995233
$ws1.Range("G14").Value = 2428635
$ws1.Range("H14").Value = -0.5902089033551768
$ws1.Range("I14").Value = "IGL"
$ws1.Range("A15").Value = "HUDCO"
$ws1.Range("B15").Value = 211
$ws1.Range("C15").Value = 211.93
$ws1.Range("D15").Value = 206.8
$ws1.Range("E15").Value = 207.25
$ws1.Range("F15").Value = 4258867
$ws1.Range("G15").Value = 10067801
$ws1.Range("H15").Value = -0.5769814083532243
$ws1.Range("I15").Value = "HUDCO"
$ws1.Range("A16").Value = "GMRAIRPORT"
$ws1.Range("B16").Value = 90.59999999999999
$ws1.Range("C16").Value = 91.09
$ws1.Range("D16").Value = 89
$ws1.Range("E16").Value = 89
$ws1.Range("F16").Value = 5554641
$ws1.Range("G16").Value = 11333128
$ws1.Range("H16").Value = -0.5098757377486604
$ws1.Range("I16").Value = "GMRAIRPORT"
$ws1.Range("A17").Value = "SYNGENE"
$ws1.Range("B17").Value = 672.15
$ws1.Range("C17").Value = 674.7
$ws1.Range("D17").Value = 651.9
$ws1.Range("E17").Value = 653.8
$ws1.Range("F17").Value = 428019
$ws1.Range("G17").Value = 953034
$ws1.Range("H17").Value = -0.5508880060942212
$ws1.Range("I17").Value = "SYNGENE"
$ws1.Range("A18").Value = "TATACHEM"
$ws1.Range("B18").Value = 959
$ws1.Range("C18").Value = 963.85
$ws1.Range("D18").Value = 941.6
$ws1.Range("E18").Value = 944.75
$ws1.Range("F18").Value = 244860
$ws1.Range("G18").Value = 586894
$ws1.Range("H18").Value = -0.5827866701653109
$ws1.Range("I18").Value = "TATACHEM"

# Sheet2 (Pos_Change) updates
$ws2.Range("A2").Value = "HDFCLIFE"
$ws2.Range("B2").Value = 756
$ws2.Range("C2").Value = 764.5
$ws2.Range("D2").Value = 756
$ws2.Range("E2").Value = 759.8
$ws2.Range("F2").Value = 2210684
$ws2.Range("G2").Value = 1423219
$ws2.Range("H2").Value = 0.5532985436535066
$ws2.Range("I2").Value = "HDFCLIFE"
$ws2.Range("A3").Value = "EICHERMOT"
$ws2.Range("B3").Value = 5681.5
$ws2.Range("C3").Value = 5711.5
$ws2.Range("D3").Value = 5655.5
$ws2.Range("E3").Value = 5670.5
$ws2.Range("F3").Value = 373874
$ws2.Range("G3").Value = 266900
$ws2.Range("H3").Value = 0.4008017984263769
$ws2.Range("I3").Value = "EICHERMOT"
$ws2.Range("A4").Value = "TATACONSUM"
$ws2.Range("B4").Value = 1052
$ws2.Range("C4").Value = 1062
$ws2.Range("D4").Value = 1043
$ws2.Range("E4").Value = 1047.9
$ws2.Range("F4").Value = 843741
$ws2.Range("G4").Value = 600790
$ws2.Range("H4").Value = 0.4043858919089865
$ws2.Range("I4").Value = "TATACONSUM"
$ws2.Range("A5").Value = "MARUTI"
$ws2.Range("B5").Value = 12630
$ws2.Range("C5").Value = 12725
$ws2.Range("D5").Value = 12564
$ws2.Range("E5").Value = 12567
$ws2.Range("F5").Value = 247048
$ws2.Range("G5").Value = 173393
$ws2.Range("H5").Value = 0.424786467735145
$ws2.Range("I5").Value = "MARUTI"
$ws2.Range("A6").Value = "WIPRO"
$ws2.Range("B6").Value = 242.4
$ws2.Range("C6").Value = 243.25
$ws2.Range("D6").Value = 238.4
$ws2.Range("E6").Value = 239.33
$ws2.Range("F6").Value = 8922847
$ws2.Range("G6").Value = 6190830
$ws2.Range("H6").Value = 0.441300601050263
$ws2.Range("I6").Value = "WIPRO"
$ws2.Range("A7").Value = "SHRIRAMFIN"
$ws2.Range("B7").Value = 623.3
$ws2.Range("C7").Value = 626.85
$ws2.Range("D7").Value = 607.95
$ws2.Range("E7").Value = 609
$ws2.Range("F7").Value = 5493423
$ws2.Range("G7").Value = 3482599
$ws2.Range("H7").Value = 0.5773917697673491
$ws2.Range("I7").Value = "SHRIRAMFIN"
$ws2.Range("A8").Value = "NAUKRI"
$ws2.Range("B8").Value = 1361
$ws2.Range("C8").Value = 1369.5
$ws2.Range("D8").Value = 1322.3
$ws2.Range("E8").Value = 1340
$ws2.Range("F8").Value = 1636607
$ws2.Range("G8").Value = 1124107
$ws2.Range("H8").Value = 0.4559174526980083
$ws2.Range("I8").Value = "NAUKRI"
$ws2.Range("A9").Value = "MPHASIS"
$ws2.Range("B9").Value = 2699
$ws2.Range("C9").Value = 2704
$ws2.Range("D9").Value = 2653.1
$ws2.Range("E9").Value = 2690.9
$ws2.Range("F9").Value = 302616
$ws2.Range("G9").Value = 200594
$ws2.Range("H9").Value = 0.5085994596049732
$ws2.Range("I9").Value = "MPHASIS"
$ws2.Range("A10").Value = "PAYTM"
$ws2.Range("B10").Value = 1065
$ws2.Range("C10").Value = 1079.9
$ws2.Range("D10").Value = 1057.6
$ws2.Range("E10").Value = 1059.9
$ws2.Range("F10").Value = 7710469
$ws2.Range("G10").Value = 5323456
$ws2.Range("H10").Value = 0.4483953657173085
$ws2.Range("I10").Value = "PAYTM"
$ws2.Range("A11").Value = "DALBHARAT"
$ws2.Range("B11").Value = 2270.7
$ws2.Range("C11").Value = 2276.8
$ws2.Range("D11").Value = 2238.2
$ws2.Range("E11").Value = 2245
$ws2.Range("F11").Value = 349347
$ws2.Range("G11").Value = 249358
$ws2.Range("H11").Value = 0.4009857313581277
$ws2.Range("I11").Value = "DALBHARAT"
$ws2.Range("A12").Value = "BSE"
$ws2.Range("B12").Value = 2500
$ws2.Range("C12").Value = 2503.6
$ws2.Range("D12").Value = 2382.3
$ws2.Range("E12").Value = 2389.1
$ws2.Range("F12").Value = 6685163
$ws2.Range("G12").Value = 4564283
$ws2.Range("H12").Value = 0.4646688209298153
$ws2.Range("I12").Value = "BSE"
$ws2.Range("A13").Value = "NYKAA"
$ws2.Range("B13").Value = 209
$ws2.Range("C13").Value = 209
$ws2.Range("D13").Value = 200.71
$ws2.Range("E13").Value = 202.5
$ws2.Range("F13").Value = 3609865
$ws2.Range("G13").Value = 2278210
$ws2.Range("H13").Value = 0.5845181085150184
$ws2.Range("I13").Value = "NYKAA"
$ws2.Range("A14").Value = "BDL"
$ws2.Range("B14").Value = 1562
$ws2.Range("C14").Value = 1567.2
$ws2.Range("D14").Value = 1491.1
$ws2.Range("E14").Value = 1493
$ws2.Range("F14").Value = 1482093
$ws2.Range("G14").Value = 993222
$ws2.Range("H14").Value = 0.4922071802678555
$ws2.Range("I14").Value = "BDL"
$ws2.Range("A15").Value = "DIXON"
$ws2.Range("B15").Value = 16690
$ws2.Range("C15").Value = 16700
$ws2.Range("D15").Value = 15810
$ws2.Range("E15").Value = 15816
$ws2.Range("F15").Value = 328403
$ws2.Range("G15").Value = 232821
$ws2.Range("H15").Value = 0.4105385682562999
$ws2.Range("I15").Value = "DIXON"
$ws2.Range("A16").Value = "KAYNES"
$ws2.Range("B16").Value = 6090
$ws2.Range("C16").Value = 6110
$ws2.Range("D16").Value = 5792.5
$ws2.Range("E16").Value = 5800.5
$ws2.Range("F16").Value = 535725
$ws2.Range("G16").Value = 355653
$ws2.Range("H16").Value = 0.5063137383910722
$ws2.Range("I16").Value = "KAYNES"

Write-Host "Done"
